# edit.ps1 - applies the diff:
# 1. Insert a new "Meta description" paragraph right after the title paragraph.
# 2. Remove the duplicated bold "Play Beast of Wealth Slot..." paragraph near the end.
# 3. Replace the final italic paragraph's text with the new image-generation prompt.

$d = $word.ActiveDocument

# --- Step 1: insert the "Meta description" paragraph after paragraph 1 ---
$titlePara = $d.Paragraphs.First
# A collapsed range positioned exactly at a paragraph boundary makes the
# runtime's InsertXML replace that paragraph instead of inserting next to it,
# so we target one character before the end (still inside paragraph 1's own
# text run) -- the new paragraph still lands immediately after paragraph 1.
$endPos = $titlePara.Range.End
$insertPoint = $d.Range($endPos - 1, $endPos - 1)

$metaSnippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our review of Beast of Wealth by Play''N''Go. Play this slot for free and enjoy exciting features, including jackpots and an RTP of 96.17%</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($metaSnippet)

# --- Step 2: delete the duplicated bold heading paragraph near the end ---
$boldParaText = "Play Beast of Wealth Slot - Free Play Included"
$count = $d.Paragraphs.Count
$dupPara = $d.Paragraphs($count - 1)
if ($dupPara.Range.Text.TrimEnd("`r") -eq $boldParaText) {
    $dupPara.Range.Delete()
}

# --- Step 3: replace the trailing italic paragraph's text with the prompt ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$oldText = "Read our review of Beast of Wealth by Play'N'Go. Play this slot for free and enjoy exciting features, including jackpots and an RTP of 96.17%"
$newText = 'Prompt: Create a cartoon-style image featuring a happy Maya warrior with glasses for the game "Beast of Wealth." The image should be visually striking and convey the excitement and adventure of the game. It should prominently feature the Maya warrior, who can be armed with various weapons and have a confident and joyful expression on his face. The background should incorporate elements of both oriental and Mayan cultures, with lush greenery, ancient temples, and bright colors. Incorporate the game logo, as well as symbols from the game, such as the dragon, tiger, turtle, and phoenix, in creative ways. Provide a sense of movement and action in the image, suggesting the energy and thrill of playing the game. Overall, the image should be both eye-catching and representative of the game''s theme and features.'
# Replace only the old-text span (not the whole paragraph range, which would
# also include the trailing paragraph mark) via direct Range.Text assignment
# -- unlike Find.Execute this keeps straight quotes/apostrophes verbatim and
# preserves the run's existing formatting (the italic rPr) and the leading
# empty run.
$startPos = $lastPara.Range.Start
$targetRange = $d.Range($startPos, $startPos + $oldText.Length)
if ($targetRange.Text -eq $oldText) {
    $targetRange.Text = $newText
}

Write-Output "done"
